$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$__s = $ws.Cells.Item(2, 4).Style
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '31.166.38'
$ws.Cells.Item(2, 4).Style = $__s
$__s = $ws.Cells.Item(2, 5).Style
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  +1.91%  '
$ws.Cells.Item(2, 5).Style = $__s

$__s = $ws.Cells.Item(3, 4).Style
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.989.99'
$ws.Cells.Item(3, 4).Style = $__s
$__s = $ws.Cells.Item(3, 5).Style
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  +5.58%  '
$ws.Cells.Item(3, 5).Style = $__s

$__s = $ws.Cells.Item(4, 4).Style
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9985'
$ws.Cells.Item(4, 4).Style = $__s
$__s = $ws.Cells.Item(4, 5).Style
$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  -0.14%  '
$ws.Cells.Item(4, 5).Style = $__s

$__s = $ws.Cells.Item(5, 4).Style
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '0.7940'
$ws.Cells.Item(5, 4).Style = $__s
$__s = $ws.Cells.Item(5, 5).Style
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  +67.58%  '
$ws.Cells.Item(5, 5).Style = $__s

$__s = $ws.Cells.Item(6, 4).Style
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '255.08'
$ws.Cells.Item(6, 4).Style = $__s
$__s = $ws.Cells.Item(6, 5).Style
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  +3.47%  '
$ws.Cells.Item(6, 5).Style = $__s

$__s = $ws.Cells.Item(7, 4).Style
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.9987'
$ws.Cells.Item(7, 4).Style = $__s
$__s = $ws.Cells.Item(7, 5).Style
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  -0.15%  '
$ws.Cells.Item(7, 5).Style = $__s

$__s = $ws.Cells.Item(8, 4).Style
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.3500'
$ws.Cells.Item(8, 4).Style = $__s
$__s = $ws.Cells.Item(8, 5).Style
$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  +20.87%  '
$ws.Cells.Item(8, 5).Style = $__s

$__s = $ws.Cells.Item(9, 4).Style
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '28.16'
$ws.Cells.Item(9, 4).Style = $__s
$__s = $ws.Cells.Item(9, 5).Style
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  +26.06%  '
$ws.Cells.Item(9, 5).Style = $__s

$__s = $ws.Cells.Item(10, 2).Style
$ws.Cells.Item(10, 2).NumberFormat = '@'
$ws.Cells.Item(10, 2).Value = 'Dogecoin'
$ws.Cells.Item(10, 2).Style = $__s
$__s = $ws.Cells.Item(10, 3).Style
$ws.Cells.Item(10, 3).NumberFormat = '@'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(10, 3).Style = $__s
$__s = $ws.Cells.Item(10, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.06995'
$ws.Cells.Item(10, 4).Style = $__s
$__s = $ws.Cells.Item(10, 5).Style
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +6.95%  '
$ws.Cells.Item(10, 5).Style = $__s

$__s = $ws.Cells.Item(11, 2).Style
$ws.Cells.Item(11, 2).NumberFormat = '@'
$ws.Cells.Item(11, 2).Value = 'Polygon'
$ws.Cells.Item(11, 2).Style = $__s
$__s = $ws.Cells.Item(11, 3).Style
$ws.Cells.Item(11, 3).NumberFormat = '@'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(11, 3).Style = $__s
$__s = $ws.Cells.Item(11, 4).Style
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.8459'
$ws.Cells.Item(11, 4).Style = $__s
$__s = $ws.Cells.Item(11, 5).Style
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  +8.96%  '
$ws.Cells.Item(11, 5).Style = $__s

$__s = $ws.Cells.Item(12, 2).Style
$ws.Cells.Item(12, 2).NumberFormat = '@'
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 2).Style = $__s
$__s = $ws.Cells.Item(12, 3).Style
$ws.Cells.Item(12, 3).NumberFormat = '@'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 3).Style = $__s
$__s = $ws.Cells.Item(12, 4).Style
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.08187'
$ws.Cells.Item(12, 4).Style = $__s
$__s = $ws.Cells.Item(12, 5).Style
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  +4.96%  '
$ws.Cells.Item(12, 5).Style = $__s

$__s = $ws.Cells.Item(13, 2).Style
$ws.Cells.Item(13, 2).NumberFormat = '@'
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 2).Style = $__s
$__s = $ws.Cells.Item(13, 3).Style
$ws.Cells.Item(13, 3).NumberFormat = '@'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 3).Style = $__s
$__s = $ws.Cells.Item(13, 4).Style
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '1.990.91'
$ws.Cells.Item(13, 4).Style = $__s
$__s = $ws.Cells.Item(13, 5).Style
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  +5.69%  '
$ws.Cells.Item(13, 5).Style = $__s

$__s = $ws.Cells.Item(14, 2).Style
$ws.Cells.Item(14, 2).NumberFormat = '@'
$ws.Cells.Item(14, 2).Value = 'Litecoin'
$ws.Cells.Item(14, 2).Style = $__s
$__s = $ws.Cells.Item(14, 3).Style
$ws.Cells.Item(14, 3).NumberFormat = '@'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(14, 3).Style = $__s
$__s = $ws.Cells.Item(14, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '100.47'
$ws.Cells.Item(14, 4).Style = $__s
$__s = $ws.Cells.Item(14, 5).Style
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  -0.58%  '
$ws.Cells.Item(14, 5).Style = $__s

$__s = $ws.Cells.Item(15, 2).Style
$ws.Cells.Item(15, 2).NumberFormat = '@'
$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 2).Style = $__s
$__s = $ws.Cells.Item(15, 3).Style
$ws.Cells.Item(15, 3).NumberFormat = '@'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(15, 3).Style = $__s
$__s = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '5.624'
$ws.Cells.Item(15, 4).Style = $__s
$__s = $ws.Cells.Item(15, 5).Style
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  +6.88%  '
$ws.Cells.Item(15, 5).Style = $__s

$__s = $ws.Cells.Item(16, 2).Style
$ws.Cells.Item(16, 2).NumberFormat = '@'
$ws.Cells.Item(16, 2).Value = 'Avalanche'
$ws.Cells.Item(16, 2).Style = $__s
$__s = $ws.Cells.Item(16, 3).Style
$ws.Cells.Item(16, 3).NumberFormat = '@'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(16, 3).Style = $__s
$__s = $ws.Cells.Item(16, 4).Style
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '15.42'
$ws.Cells.Item(16, 4).Style = $__s
$__s = $ws.Cells.Item(16, 5).Style
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  +16.66%  '
$ws.Cells.Item(16, 5).Style = $__s

$__s = $ws.Cells.Item(17, 2).Style
$ws.Cells.Item(17, 2).NumberFormat = '@'
$ws.Cells.Item(17, 2).Value = 'BitcoinCash'
$ws.Cells.Item(17, 2).Style = $__s
$__s = $ws.Cells.Item(17, 3).Style
$ws.Cells.Item(17, 3).NumberFormat = '@'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(17, 3).Style = $__s
$__s = $ws.Cells.Item(17, 4).Style
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '273.51'
$ws.Cells.Item(17, 4).Style = $__s
$__s = $ws.Cells.Item(17, 5).Style
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  -4.14%  '
$ws.Cells.Item(17, 5).Style = $__s

$__s = $ws.Cells.Item(18, 2).Style
$ws.Cells.Item(18, 2).NumberFormat = '@'
$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 2).Style = $__s
$__s = $ws.Cells.Item(18, 3).Style
$ws.Cells.Item(18, 3).NumberFormat = '@'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 3).Style = $__s
$__s = $ws.Cells.Item(18, 4).Style
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '31.161.67'
$ws.Cells.Item(18, 4).Style = $__s
$__s = $ws.Cells.Item(18, 5).Style
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  +1.98%  '
$ws.Cells.Item(18, 5).Style = $__s

$__s = $ws.Cells.Item(19, 2).Style
$ws.Cells.Item(19, 2).NumberFormat = '@'
$ws.Cells.Item(19, 2).Value = 'Uniswap'
$ws.Cells.Item(19, 2).Style = $__s
$__s = $ws.Cells.Item(19, 3).Style
$ws.Cells.Item(19, 3).NumberFormat = '@'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(19, 3).Style = $__s
$__s = $ws.Cells.Item(19, 4).Style
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '5.871'
$ws.Cells.Item(19, 4).Style = $__s
$__s = $ws.Cells.Item(19, 5).Style
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  +9.72%  '
$ws.Cells.Item(19, 5).Style = $__s

$__s = $ws.Cells.Item(20, 2).Style
$ws.Cells.Item(20, 2).NumberFormat = '@'
$ws.Cells.Item(20, 2).Value = 'ShibaInu'
$ws.Cells.Item(20, 2).Style = $__s
$__s = $ws.Cells.Item(20, 3).Style
$ws.Cells.Item(20, 3).NumberFormat = '@'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(20, 3).Style = $__s
$__s = $ws.Cells.Item(20, 4).Style
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '0.000007928'
$ws.Cells.Item(20, 4).Style = $__s
$__s = $ws.Cells.Item(20, 5).Style
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  +5.28%  '
$ws.Cells.Item(20, 5).Style = $__s

$__s = $ws.Cells.Item(21, 2).Style
$ws.Cells.Item(21, 2).NumberFormat = '@'
$ws.Cells.Item(21, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(21, 2).Style = $__s
$__s = $ws.Cells.Item(21, 3).Style
$ws.Cells.Item(21, 3).NumberFormat = '@'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(21, 3).Style = $__s
$__s = $ws.Cells.Item(21, 4).Style
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '2.252.53'
$ws.Cells.Item(21, 4).Style = $__s
$__s = $ws.Cells.Item(21, 5).Style
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  +5.95%  '
$ws.Cells.Item(21, 5).Style = $__s

$__s = $ws.Cells.Item(22, 2).Style
$ws.Cells.Item(22, 2).NumberFormat = '@'
$ws.Cells.Item(22, 2).Value = 'Dai'
$ws.Cells.Item(22, 2).Style = $__s
$__s = $ws.Cells.Item(22, 3).Style
$ws.Cells.Item(22, 3).NumberFormat = '@'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(22, 3).Style = $__s
$__s = $ws.Cells.Item(22, 4).Style
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.9990'
$ws.Cells.Item(22, 4).Style = $__s
$__s = $ws.Cells.Item(22, 5).Style
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  -0.14%  '
$ws.Cells.Item(22, 5).Style = $__s

$__s = $ws.Cells.Item(23, 2).Style
$ws.Cells.Item(23, 2).NumberFormat = '@'
$ws.Cells.Item(23, 2).Value = 'BinanceUSD'
$ws.Cells.Item(23, 2).Style = $__s
$__s = $ws.Cells.Item(23, 3).Style
$ws.Cells.Item(23, 3).NumberFormat = '@'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(23, 3).Style = $__s
$__s = $ws.Cells.Item(23, 4).Style
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.9991'
$ws.Cells.Item(23, 4).Style = $__s
$__s = $ws.Cells.Item(23, 5).Style
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  -0.12%  '
$ws.Cells.Item(23, 5).Style = $__s

$__s = $ws.Cells.Item(24, 2).Style
$ws.Cells.Item(24, 2).NumberFormat = '@'
$ws.Cells.Item(24, 2).Value = 'Chainlink'
$ws.Cells.Item(24, 2).Style = $__s
$__s = $ws.Cells.Item(24, 3).Style
$ws.Cells.Item(24, 3).NumberFormat = '@'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(24, 3).Style = $__s
$__s = $ws.Cells.Item(24, 4).Style
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '7.056'
$ws.Cells.Item(24, 4).Style = $__s
$__s = $ws.Cells.Item(24, 5).Style
$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  +9.75%  '
$ws.Cells.Item(24, 5).Style = $__s

$__s = $ws.Cells.Item(25, 2).Style
$ws.Cells.Item(25, 2).NumberFormat = '@'
$ws.Cells.Item(25, 2).Value = 'Cosmos'
$ws.Cells.Item(25, 2).Style = $__s
$__s = $ws.Cells.Item(25, 3).Style
$ws.Cells.Item(25, 3).NumberFormat = '@'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(25, 3).Style = $__s
$__s = $ws.Cells.Item(25, 4).Style
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '9.998'
$ws.Cells.Item(25, 4).Style = $__s
$__s = $ws.Cells.Item(25, 5).Style
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  +9.11%  '
$ws.Cells.Item(25, 5).Style = $__s

$__s = $ws.Cells.Item(26, 2).Style
$ws.Cells.Item(26, 2).NumberFormat = '@'
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 2).Style = $__s
$__s = $ws.Cells.Item(26, 3).Style
$ws.Cells.Item(26, 3).NumberFormat = '@'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(26, 3).Style = $__s
$__s = $ws.Cells.Item(26, 4).Style
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '165.76'
$ws.Cells.Item(26, 4).Style = $__s
$__s = $ws.Cells.Item(26, 5).Style
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  +2.13%  '
$ws.Cells.Item(26, 5).Style = $__s

$__s = $ws.Cells.Item(27, 2).Style
$ws.Cells.Item(27, 2).NumberFormat = '@'
$ws.Cells.Item(27, 2).Value = 'Stellar'
$ws.Cells.Item(27, 2).Style = $__s
$__s = $ws.Cells.Item(27, 3).Style
$ws.Cells.Item(27, 3).NumberFormat = '@'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(27, 3).Style = $__s
$__s = $ws.Cells.Item(27, 4).Style
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '0.1506'
$ws.Cells.Item(27, 4).Style = $__s
$__s = $ws.Cells.Item(27, 5).Style
$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  +55.20%  '
$ws.Cells.Item(27, 5).Style = $__s

$__s = $ws.Cells.Item(28, 2).Style
$ws.Cells.Item(28, 2).NumberFormat = '@'
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 2).Style = $__s
$__s = $ws.Cells.Item(28, 3).Style
$ws.Cells.Item(28, 3).NumberFormat = '@'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28, 3).Style = $__s
$__s = $ws.Cells.Item(28, 4).Style
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '19.92'
$ws.Cells.Item(28, 4).Style = $__s
$__s = $ws.Cells.Item(28, 5).Style
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  +4.13%  '
$ws.Cells.Item(28, 5).Style = $__s

$__s = $ws.Cells.Item(29, 2).Style
$ws.Cells.Item(29, 2).NumberFormat = '@'
$ws.Cells.Item(29, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(29, 2).Style = $__s
$__s = $ws.Cells.Item(29, 3).Style
$ws.Cells.Item(29, 3).NumberFormat = '@'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(29, 3).Style = $__s
$__s = $ws.Cells.Item(29, 4).Style
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.346'
$ws.Cells.Item(29, 4).Style = $__s
$__s = $ws.Cells.Item(29, 5).Style
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  +22.42%  '
$ws.Cells.Item(29, 5).Style = $__s

$__s = $ws.Cells.Item(30, 2).Style
$ws.Cells.Item(30, 2).NumberFormat = '@'
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 2).Style = $__s
$__s = $ws.Cells.Item(30, 3).Style
$ws.Cells.Item(30, 3).NumberFormat = '@'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 3).Style = $__s
$__s = $ws.Cells.Item(30, 4).Style
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.597'
$ws.Cells.Item(30, 4).Style = $__s
$__s = $ws.Cells.Item(30, 5).Style
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  +6.26%  '
$ws.Cells.Item(30, 5).Style = $__s

$__s = $ws.Cells.Item(31, 2).Style
$ws.Cells.Item(31, 2).NumberFormat = '@'
$ws.Cells.Item(31, 2).Value = 'Toncoin'
$ws.Cells.Item(31, 2).Style = $__s
$__s = $ws.Cells.Item(31, 3).Style
$ws.Cells.Item(31, 3).NumberFormat = '@'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(31, 3).Style = $__s
$__s = $ws.Cells.Item(31, 4).Style
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.357'
$ws.Cells.Item(31, 4).Style = $__s
$__s = $ws.Cells.Item(31, 5).Style
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  +1.87%  '
$ws.Cells.Item(31, 5).Style = $__s

$__s = $ws.Cells.Item(32, 2).Style
$ws.Cells.Item(32, 2).NumberFormat = '@'
$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 2).Style = $__s
$__s = $ws.Cells.Item(32, 3).Style
$ws.Cells.Item(32, 3).NumberFormat = '@'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(32, 3).Style = $__s
$__s = $ws.Cells.Item(32, 4).Style
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '4.589'
$ws.Cells.Item(32, 4).Style = $__s
$__s = $ws.Cells.Item(32, 5).Style
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  +7.63%  '
$ws.Cells.Item(32, 5).Style = $__s

$__s = $ws.Cells.Item(33, 2).Style
$ws.Cells.Item(33, 2).NumberFormat = '@'
$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 2).Style = $__s
$__s = $ws.Cells.Item(33, 3).Style
$ws.Cells.Item(33, 3).NumberFormat = '@'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(33, 3).Style = $__s
$__s = $ws.Cells.Item(33, 4).Style
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '4.405'
$ws.Cells.Item(33, 4).Style = $__s
$__s = $ws.Cells.Item(33, 5).Style
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  +5.01%  '
$ws.Cells.Item(33, 5).Style = $__s

$__s = $ws.Cells.Item(34, 2).Style
$ws.Cells.Item(34, 2).NumberFormat = '@'
$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 2).Style = $__s
$__s = $ws.Cells.Item(34, 3).Style
$ws.Cells.Item(34, 3).NumberFormat = '@'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(34, 3).Style = $__s
$__s = $ws.Cells.Item(34, 4).Style
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.05244'
$ws.Cells.Item(34, 4).Style = $__s
$__s = $ws.Cells.Item(34, 5).Style
$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  +8.33%  '
$ws.Cells.Item(34, 5).Style = $__s

$__s = $ws.Cells.Item(35, 2).Style
$ws.Cells.Item(35, 2).NumberFormat = '@'
$ws.Cells.Item(35, 2).Value = 'ImmutableX'
$ws.Cells.Item(35, 2).Style = $__s
$__s = $ws.Cells.Item(35, 3).Style
$ws.Cells.Item(35, 3).NumberFormat = '@'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(35, 3).Style = $__s
$__s = $ws.Cells.Item(35, 4).Style
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.7799'
$ws.Cells.Item(35, 4).Style = $__s
$__s = $ws.Cells.Item(35, 5).Style
$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  +11.82%  '
$ws.Cells.Item(35, 5).Style = $__s

$__s = $ws.Cells.Item(36, 4).Style
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '1.221'
$ws.Cells.Item(36, 4).Style = $__s
$__s = $ws.Cells.Item(36, 5).Style
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  +8.01%  '
$ws.Cells.Item(36, 5).Style = $__s

$__s = $ws.Cells.Item(37, 2).Style
$ws.Cells.Item(37, 2).NumberFormat = '@'
$ws.Cells.Item(37, 2).Value = 'HuobiToken'
$ws.Cells.Item(37, 2).Style = $__s
$__s = $ws.Cells.Item(37, 3).Style
$ws.Cells.Item(37, 3).NumberFormat = '@'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(37, 3).Style = $__s
$__s = $ws.Cells.Item(37, 4).Style
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.764'
$ws.Cells.Item(37, 4).Style = $__s
$__s = $ws.Cells.Item(37, 5).Style
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  +0.27%  '
$ws.Cells.Item(37, 5).Style = $__s

$__s = $ws.Cells.Item(38, 2).Style
$ws.Cells.Item(38, 2).NumberFormat = '@'
$ws.Cells.Item(38, 2).Value = 'Frax'
$ws.Cells.Item(38, 2).Style = $__s
$__s = $ws.Cells.Item(38, 3).Style
$ws.Cells.Item(38, 3).NumberFormat = '@'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(38, 3).Style = $__s
$__s = $ws.Cells.Item(38, 4).Style
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.9982'
$ws.Cells.Item(38, 4).Style = $__s
$__s = $ws.Cells.Item(38, 5).Style
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  -0.17%  '
$ws.Cells.Item(38, 5).Style = $__s

$__s = $ws.Cells.Item(39, 4).Style
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.01999'
$ws.Cells.Item(39, 4).Style = $__s
$__s = $ws.Cells.Item(39, 5).Style
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  +4.56%  '
$ws.Cells.Item(39, 5).Style = $__s

$__s = $ws.Cells.Item(40, 4).Style
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.878'
$ws.Cells.Item(40, 4).Style = $__s
$__s = $ws.Cells.Item(40, 5).Style
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  -0.64%  '
$ws.Cells.Item(40, 5).Style = $__s

$__s = $ws.Cells.Item(41, 4).Style
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '6.647'
$ws.Cells.Item(41, 4).Style = $__s
$__s = $ws.Cells.Item(41, 5).Style
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  +5.47%  '
$ws.Cells.Item(41, 5).Style = $__s

$__s = $ws.Cells.Item(42, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '79.68'
$ws.Cells.Item(42, 4).Style = $__s
$__s = $ws.Cells.Item(42, 5).Style
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  +4.81%  '
$ws.Cells.Item(42, 5).Style = $__s

$__s = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.4675'
$ws.Cells.Item(43, 4).Style = $__s
$__s = $ws.Cells.Item(43, 5).Style
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  +9.68%  '
$ws.Cells.Item(43, 5).Style = $__s

$__s = $ws.Cells.Item(44, 4).Style
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.121'
$ws.Cells.Item(44, 4).Style = $__s
$__s = $ws.Cells.Item(44, 5).Style
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  +6.91%  '
$ws.Cells.Item(44, 5).Style = $__s

$__s = $ws.Cells.Item(45, 4).Style
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.8504'
$ws.Cells.Item(45, 4).Style = $__s
$__s = $ws.Cells.Item(45, 5).Style
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  +2.32%  '
$ws.Cells.Item(45, 5).Style = $__s

$__s = $ws.Cells.Item(46, 4).Style
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '104.51'
$ws.Cells.Item(46, 4).Style = $__s
$__s = $ws.Cells.Item(46, 5).Style
$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  +2.92%  '
$ws.Cells.Item(46, 5).Style = $__s

$__s = $ws.Cells.Item(47, 4).Style
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.9989'
$ws.Cells.Item(47, 4).Style = $__s
$__s = $ws.Cells.Item(47, 5).Style
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  -0.11%  '
$ws.Cells.Item(47, 5).Style = $__s

$__s = $ws.Cells.Item(48, 4).Style
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '7.691'
$ws.Cells.Item(48, 4).Style = $__s
$__s = $ws.Cells.Item(48, 5).Style
$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  +9.40%  '
$ws.Cells.Item(48, 5).Style = $__s

$__s = $ws.Cells.Item(49, 4).Style
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '9.888'
$ws.Cells.Item(49, 4).Style = $__s
$__s = $ws.Cells.Item(49, 5).Style
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  +0.37%  '
$ws.Cells.Item(49, 5).Style = $__s

$__s = $ws.Cells.Item(50, 4).Style
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '36.79'
$ws.Cells.Item(50, 4).Style = $__s
$__s = $ws.Cells.Item(50, 5).Style
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  +4.42%  '
$ws.Cells.Item(50, 5).Style = $__s

$__s = $ws.Cells.Item(51, 5).Style
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  +8.48%  '
$ws.Cells.Item(51, 5).Style = $__s
